# Update income statement worksheet:
#  - add a new "units" header column (H) formatted like the other header cells
#  - add empty text placeholder cells in H2:H14
#  - revise several figures in columns C:G for rows 6, 8, 9, 10, 12, 14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column H -------------------------------------------------------
# Give H1 the same (bold/centered/bordered) header style as the existing
# header cells by copying G1's formatting, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "units"

# H2:H14 become empty (but present) text cells. A leading apostrophe makes
# the runtime store a literal/text empty value instead of clearing the
# cell entirely.
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 8).Value = "'"
}

# --- Revised figures ------------------------------------------------------
# Row 6 - Selling, General & Administrative
$ws.Range("E6").Value = "'-257"
$ws.Range("F6").Value = "'-258"
$ws.Range("G6").Value = "'-269"

# Row 8 - Total Operating Expenses
$ws.Range("C8").Value = "'330"
$ws.Range("D8").Value = "'405"
$ws.Range("E8").Value = "'423"
$ws.Range("F8").Value = "'432"
$ws.Range("G8").Value = "'466"

# Row 9 - EBITDA
$ws.Range("C9").Value = "'1,700"
$ws.Range("D9").Value = "'2,070"
$ws.Range("E9").Value = "'2,332"
$ws.Range("F9").Value = "'2,628"
$ws.Range("G9").Value = "'3,001"

# Row 10 - Operating Income (EBIT)
$ws.Range("E10").Value = "'1,522"
$ws.Range("F10").Value = "'1,801"
$ws.Range("G10").Value = "'2,108"

# Row 12 - Earnings Before Tax
$ws.Range("C12").Value = "'1,120"
$ws.Range("D12").Value = "'1,355"
$ws.Range("E12").Value = "'1,585"
$ws.Range("F12").Value = "'1,867"
$ws.Range("G12").Value = "'2,178"

# Row 14 - Net Income
$ws.Range("C14").Value = "'1,195"
$ws.Range("D14").Value = "'1,445"
$ws.Range("E14").Value = "'1,693"
$ws.Range("F14").Value = "'1,997"
$ws.Range("G14").Value = "'2,333"
